# Update the "想去人数" (want-to-go count) figures in column F for rows 2-7
# on both the "展览" and "全部类型" worksheets (they carry duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1364
    3 = 2096
    4 = 280
    5 = 73
    6 = 6378
    7 = 261
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
